# Update TPM-derived LR-pair statistics (Fgf9-Fgfr3) to reflect new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.520102666666667
$ws.Range("H2").Value = 13.560308
$ws.Range("I2").Value = 0.9927775608668273
$ws.Range("J2").Value = 0.9927775608668273
$ws.Range("M2").Value = 4.959409333333333
$ws.Range("N2").Value = 14.878228
$ws.Range("O2").Value = 0.8271666313262851
$ws.Range("P2").Value = 0.8271666313262852
$ws.Range("Q2").Value = 22.41703935269156
$ws.Range("R2").Value = 201.753354174224
$ws.Range("S2").Value = 0.8211924706785395
$ws.Range("T2").Value = 0.8211924706785396

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.520102666666667
$ws.Range("H3").Value = 13.560308
$ws.Range("I3").Value = 0.9927775608668273
$ws.Range("J3").Value = 0.9927775608668273
$ws.Range("O3").Value = 0.09421438109281059
$ws.Range("P3").Value = 0.09421438109281059
$ws.Range("Q3").Value = 2.553303540739556
$ws.Range("R3").Value = 22.979731866656
$ws.Range("S3").Value = 0.09353392345989824
$ws.Range("T3").Value = 0.09353392345989824

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.520102666666667
$ws.Range("H4").Value = 13.560308
$ws.Range("I4").Value = 0.9927775608668273
$ws.Range("J4").Value = 0.9927775608668273
$ws.Range("O4").Value = 0.07861898758090437
$ws.Range("P4").Value = 0.07861898758090438
$ws.Range("Q4").Value = 2.130652847593778
$ws.Range("R4").Value = 19.175875628344
$ws.Range("S4").Value = 0.07805116672838963
$ws.Range("T4").Value = 0.07805116672838965

# Row 5
$ws.Range("I5").Value = 0.007222439133172593
$ws.Range("J5").Value = 0.007222439133172593
$ws.Range("M5").Value = 4.959409333333333
$ws.Range("N5").Value = 14.878228
$ws.Range("O5").Value = 0.8271666313262851
$ws.Range("P5").Value = 0.8271666313262852
$ws.Range("Q5").Value = 0.1630835633808889
$ws.Range("R5").Value = 1.467752070428
$ws.Range("S5").Value = 0.005974160647745508
$ws.Range("T5").Value = 0.005974160647745508

# Row 6
$ws.Range("I6").Value = 0.007222439133172593
$ws.Range("J6").Value = 0.007222439133172593
$ws.Range("O6").Value = 0.09421438109281059
$ws.Range("P6").Value = 0.09421438109281059
$ws.Range("S6").Value = 0.0006804576329123512
$ws.Range("T6").Value = 0.0006804576329123512

# Row 7
$ws.Range("I7").Value = 0.007222439133172593
$ws.Range("J7").Value = 0.007222439133172593
$ws.Range("O7").Value = 0.07861898758090437
$ws.Range("P7").Value = 0.07861898758090438
$ws.Range("S7").Value = 0.0005678208525147337
$ws.Range("T7").Value = 0.0005678208525147339
